# Applies the "English words added 3.3.0" edit to the Gilaki word list.
#
# 1. Rows 11-70 (glk_word / glk_example / en_word / en_example) are rotated
#    up by one: the entry that was on row 12 moves to row 11, row 13 -> 12,
#    ... row 70 -> 69, and the entry that used to be on row 11 (the "Walnut"
#    / "آغۊز" entry) wraps around onto row 70 - with its leading alef
#    corrected from alef-madda (آ) to alef-hamza (أ).
# 2. A handful of other "Walnut" mentions elsewhere in the sheet get the
#    same alef-madda -> alef-hamza spelling correction (column B only).
# 3. A batch of example sentences that contain the superscript-alef mark
#    (ٰ) gain one extra space before that mark.
# 4. Row 611's Gilaki headword gains a second word.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rotate rows 11-70 up by one -----------------------------------
# Stash the row-11 values (they'll be overwritten by the bulk copy below).
$origA11 = $ws.Cells.Item(11, 1).Value()
$origB11 = $ws.Cells.Item(11, 2).Value()
$origC11 = $ws.Cells.Item(11, 3).Value()
$origD11 = $ws.Cells.Item(11, 4).Value()

# Bulk-copy rows 12-70 onto rows 11-69 in one shot.
$srcVals = $ws.Range("A12:D70").Value()
$ws.Range("A11:D69").Value = $srcVals

# Row 70 becomes the old row 11 ("Walnut"), with the alef spelling fixed.
$ws.Cells.Item(70, 1).Value = $origA11.Replace([string][char]0x0622, [string][char]0x0623)
$ws.Cells.Item(70, 2).Value = $origB11.Replace([string][char]0x0622, [string][char]0x0623)
$ws.Cells.Item(70, 3).Value = $origC11
$ws.Cells.Item(70, 4).Value = $origD11

# --- 2. Fix the alef-madda -> alef-hamza spelling of the walnut word --
#        in the remaining example sentences that mention it.
$madda = [string][char]0x0622
$hamza = [string][char]0x0623
$walnutRows = @(174, 206, 372, 487, 768, 773)
foreach ($r in $walnutRows) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.Value = $cell.Value().Replace($madda, $hamza)
}

# --- 3. Insert one extra space before the superscript alef (ٰ) mark ---
$supAlef = [string][char]0x0670
$spaceFixRows = @(106, 189, 209, 298, 316, 330, 359, 371, 493, 562, 586, 665, 720, 724, 751, 799, 802, 811, 848, 926, 944)
foreach ($r in $spaceFixRows) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.Value = $cell.Value().Replace($supAlef, " " + $supAlef)
}

# --- 4. Row 611 headword gains a second word ---------------------------
$ws.Cells.Item(611, 1).Value = $ws.Cells.Item(611, 1).Value() + " دأئن"
